$wb = $excel.ActiveWorkbook

# --- Update the "Conversión del día" text block on Hoja1 ---
$wsHoja1 = $wb.Worksheets.Item("Hoja1")
$newText = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 6.58 = 26118.42 pesos`n✅ 26118.42 pesos = 6.56 = 977.8 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"
$wsHoja1.Range("A1").Value = $newText

# --- Update tasas rates on "tasas" sheet ---
$wsTasas = $wb.Worksheets.Item("tasas")
$wsTasas.Range("N10").Value = 152
$wsTasas.Range("O10").Value = 3970
$wsTasas.Range("N12").Value = 3979.99
$wsTasas.Range("O12").Value = 149
